$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AJ1").Value = 0.73778547707250342
$ws.Range("BJ1").Value = 0.77762020732056381
$ws.Range("A2").Value = 0.76922673442555778
$ws.Range("BD2").Value = 0.91511651502192515
$ws.Range("S3").Value = 0.95927744509921842
$ws.Range("AD3").Value = 0.88520959997945048
$ws.Range("AO3").Value = 0.69775661291387803
$ws.Range("AT3").Value = 0.85805192766325811
$ws.Range("V4").Value = 0.96124099297667032
$ws.Range("E6").Value = 0.70266144114894413
$ws.Range("AV6").Value = 0.97691413482895073
$ws.Range("I7").Value = 0.59374976612030284
$ws.Range("AY7").Value = 0.68117467304503232
$ws.Range("Q8").Value = 0.68291577613662047
$ws.Range("AL8").Value = 0.79636947180141981
$ws.Range("AT8").Value = 0.9547663708881291
$ws.Range("BI8").Value = 0.91383404814341851
$ws.Range("K9").Value = 0.85122450899930069
$ws.Range("BI9").Value = 0.61629769599858764
$ws.Range("I10").Value = 0.9287463195003931
$ws.Range("AA11").Value = 0.95941525563344654
$ws.Range("BN11").Value = 0.99220170563480503
$ws.Range("J12").Value = 0.76952637758525078
$ws.Range("V13").Value = 0.76680126122733494
$ws.Range("BD13").Value = 0.73213463010481594
$ws.Range("BA14").Value = 0.75772862256168172
$ws.Range("BD15").Value = 0.92075273291565141
$ws.Range("BG15").Value = 0.91297618023856186
$ws.Range("BE16").Value = 0.86322339961884231
$ws.Range("AG18").Value = 0.98530639964696209
$ws.Range("AO18").Value = 0.87785548221056264
$ws.Range("BN18").Value = 0.80640574520615083
$ws.Range("R19").Value = 0.97237855132344686
$ws.Range("AU19").Value = 0.74869449764690188
$ws.Range("AR20").Value = 0.97465880742577726
$ws.Range("B21").Value = 0.93598886228106037
$ws.Range("BP21").Value = 0.90216932484964807
$ws.Range("T22").Value = 0.98168572612343619
$ws.Range("U22").Value = 0.64632499251103503
$ws.Range("Y24").Value = 0.87865380406379079
$ws.Range("AK24").Value = 0.79241197017828791
$ws.Range("BH24").Value = 0.97115459370378665
$ws.Range("Z25").Value = 0.96536971813419359
$ws.Range("BJ25").Value = 0.79636804554167173
$ws.Range("Q27").Value = 0.66330704777967098
$ws.Range("W27").Value = 0.66796361524285763
$ws.Range("AN27").Value = 0.74439070793696804
$ws.Range("E29").Value = 0.78540670711257965
$ws.Range("C31").Value = 0.85898050082293098
$ws.Range("E31").Value = 0.91313346199885592
$ws.Range("AM31").Value = 0.6657507012505377
$ws.Range("AR31").Value = 0.91935941758116768
$ws.Range("AT32").Value = 0.94938074679339379
$ws.Range("AT33").Value = 0.5561582943546014
$ws.Range("BG34").Value = 0.96089454556724507
$ws.Range("W35").Value = 0.61583237704997651
$ws.Range("X35").Value = 0.89625066628704908
$ws.Range("AG35").Value = 0.91508381350417589
$ws.Range("AH35").Value = 0.86272480660172335
$ws.Range("BN35").Value = 0.52770270195392388
$ws.Range("BH36").Value = 0.83894389962053773
$ws.Range("AB37").Value = 0.81212432460852158
$ws.Range("AT37").Value = 0.95395874667094949
$ws.Range("AX37").Value = 0.94255303732603679
$ws.Range("AY37").Value = 0.98022606607880225
$ws.Range("BD37").Value = 0.9471800760202389
$ws.Range("BI37").Value = 0.83895844340783643
$ws.Range("AJ38").Value = 0.65098135458682993
$ws.Range("L39").Value = 0.83515752144534705
$ws.Range("AY39").Value = 0.81583350961448575
$ws.Range("Q40").Value = 0.95247001563561162
$ws.Range("AF40").Value = 0.6822757053812265
$ws.Range("K41").Value = 0.74356607340947534
$ws.Range("U41").Value = 0.78872797790276572
$ws.Range("AU41").Value = 0.90010623442621296
$ws.Range("F42").Value = 0.99978144037216754
$ws.Range("AM42").Value = 0.9359095735566092
$ws.Range("AN42").Value = 0.79486927419551678
$ws.Range("AV42").Value = 0.69408148116275958
$ws.Range("AW42").Value = 0.9491662083366148
$ws.Range("BL42").Value = 0.99989785643840146
$ws.Range("AD43").Value = 0.95997108158504185
$ws.Range("AL43").Value = 0.68612965012522165
$ws.Range("BJ43").Value = 0.51737871615732278
$ws.Range("AO44").Value = 0.88483537687434732
$ws.Range("J45").Value = 0.92352861980142043
$ws.Range("D46").Value = 0.82294320587956249
$ws.Range("AM46").Value = 0.96706017094355723
$ws.Range("BJ46").Value = 0.71742113631641802
$ws.Range("M47").Value = 0.67440657258428161
$ws.Range("Q47").Value = 0.64868011805318493
$ws.Range("BF49").Value = 0.60287353733656279
$ws.Range("BN49").Value = 0.58958795899174543
$ws.Range("W50").Value = 0.65986730666184656
$ws.Range("BE50").Value = 0.62149192551100896
$ws.Range("N51").Value = 0.91545652118685961
$ws.Range("AS51").Value = 0.96463664189798992
$ws.Range("AV51").Value = 0.56563497640301286
$ws.Range("AX52").Value = 0.71164716495865488
$ws.Range("BG52").Value = 0.74587091658735027
$ws.Range("D53").Value = 0.87173266969576813
$ws.Range("P53").Value = 0.88186754718059057
$ws.Range("T53").Value = 0.91523988240003551
$ws.Range("AB53").Value = 0.5697858421380565
$ws.Range("AM53").Value = 0.91965534247466507
$ws.Range("BP53").Value = 0.83108341157279764
$ws.Range("Z54").Value = 0.95447221566449436
$ws.Range("AC54").Value = 0.87754715202564026
$ws.Range("AQ54").Value = 0.8094374670931892
$ws.Range("AZ54").Value = 0.72207781259293102
$ws.Range("BC54").Value = 0.9971420471462491
$ws.Range("Y55").Value = 0.9060372007606009
$ws.Range("AE55").Value = 0.75599175248572936
$ws.Range("N56").Value = 0.70510250211135261
$ws.Range("W56").Value = 0.82642961205927068
$ws.Range("BL57").Value = 0.84047815408948501
$ws.Range("AS58").Value = 0.97892476574148524
$ws.Range("BD58").Value = 0.78828911858405393
$ws.Range("BF59").Value = 0.98622024572136002
$ws.Range("BI59").Value = 0.9621937508804157
$ws.Range("AF62").Value = 0.81036225702015141
$ws.Range("AP62").Value = 0.84885148976848157
$ws.Range("Z63").Value = 0.96126076741806288
$ws.Range("AF63").Value = 0.8971726007610521
$ws.Range("J64").Value = 0.79382170956459075
$ws.Range("AS64").Value = 0.96049306529920453
$ws.Range("BM64").Value = 0.76816543642380042
$ws.Range("AU65").Value = 0.99362696012029184
$ws.Range("AA66").Value = 0.80504578979385699
$ws.Range("AQ66").Value = 0.90929727755710654
$ws.Range("BK66").Value = 0.69471829268454788
$ws.Range("BO66").Value = 0.91521529037906657
$ws.Range("AJ67").Value = 0.84054198130171676
$ws.Range("AR67").Value = 0.9063256062004158
$ws.Range("BM67").Value = 0.98009246000586414
$ws.Range("I68").Value = 0.81905719558184265
